$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.31"
$ws.Range("E2").Value = "'-0.31%"
$ws.Range("G2").Value = "'5"
$ws.Range("D3").Value = "'39.74"
$ws.Range("E3").Value = "'2.30%"
$ws.Range("G3").Value = "'5"
$ws.Range("D4").Value = "'5.161"
$ws.Range("E4").Value = "'1.40%"
$ws.Range("G4").Value = "'5"
$ws.Range("D5").Value = "'0.08137"
$ws.Range("E5").Value = "'-0.59%"
$ws.Range("G5").Value = "'5"
$ws.Range("D6").Value = "'1.946"
$ws.Range("E6").Value = "'-2.38%"
$ws.Range("G6").Value = "'5"
$ws.Range("D7").Value = "'8.152"
$ws.Range("E7").Value = "'3.06%"
$ws.Range("G7").Value = "'5"
$ws.Range("D8").Value = "'4.233"
$ws.Range("E8").Value = "'1.35%"
$ws.Range("G8").Value = "'5"
$ws.Range("D9").Value = "'0.9296"
$ws.Range("E9").Value = "'-0.40%"
$ws.Range("G9").Value = "'5"
$ws.Range("D10").Value = "'0.1428"
$ws.Range("E10").Value = "'1.52%"
$ws.Range("G10").Value = "'5"
$ws.Range("D11").Value = "'0.1932"
$ws.Range("E11").Value = "'-0.94%"
$ws.Range("G11").Value = "'5"
$ws.Range("D12").Value = "'0.09134"
$ws.Range("E12").Value = "'-1.64%"
$ws.Range("G12").Value = "'5"
$ws.Range("E13").Value = "'1.29%"
$ws.Range("G13").Value = "'5"
$ws.Range("D14").Value = "'0.09798"
$ws.Range("E14").Value = "'-0.45%"
$ws.Range("G14").Value = "'5"
$ws.Range("D15").Value = "'0.001396"
$ws.Range("E15").Value = "'-1.46%"
$ws.Range("G15").Value = "'5"
$ws.Range("D16").Value = "'0.005842"
$ws.Range("E16").Value = "'-4.33%"
$ws.Range("G16").Value = "'5"
$ws.Range("D17").Value = "'3.924"
$ws.Range("E17").Value = "'4.17%"
$ws.Range("G17").Value = "'5"
$ws.Range("D18").Value = "'3.322"
$ws.Range("E18").Value = "'-4.27%"
$ws.Range("G18").Value = "'5"
$ws.Range("D19").Value = "'0.3430"
$ws.Range("E19").Value = "'-0.61%"
$ws.Range("G19").Value = "'5"
$ws.Range("D20").Value = "'0.1321"
$ws.Range("E20").Value = "'1.34%"
$ws.Range("G20").Value = "'5"
$ws.Range("D21").Value = "'4.631"
$ws.Range("E21").Value = "'-3.90%"
$ws.Range("G21").Value = "'5"
$ws.Range("E22").Value = "'-0.13%"
$ws.Range("G22").Value = "'5"
$ws.Range("D23").Value = "'0.04377"
$ws.Range("E23").Value = "'-1.88%"
$ws.Range("G23").Value = "'5"
$ws.Range("E24").Value = "'-1.28%"
$ws.Range("G24").Value = "'5"
$ws.Range("D25").Value = "'0.004374"
$ws.Range("E25").Value = "'4.80%"
$ws.Range("G25").Value = "'5"
$ws.Range("D26").Value = "'0.0001301"
$ws.Range("E26").Value = "'-0.08%"
$ws.Range("G26").Value = "'5"
$ws.Range("E27").Value = "'-10.05%"
$ws.Range("G27").Value = "'5"
$ws.Range("G28").Value = "'5"
$ws.Range("G29").Value = "'5"
$ws.Range("G30").Value = "'5"
$ws.Range("G31").Value = "'5"
$ws.Range("G32").Value = "'5"
$ws.Range("G33").Value = "'5"
$ws.Range("G34").Value = "'5"
$ws.Range("G35").Value = "'5"
$ws.Range("G36").Value = "'5"
$ws.Range("G37").Value = "'5"
$ws.Range("G38").Value = "'5"
$ws.Range("D39").Value = "'0.02049"
$ws.Range("E39").Value = "'-3.84%"
$ws.Range("G39").Value = "'5"
$ws.Range("D40").Value = "'0.05078"
$ws.Range("E40").Value = "'-1.82%"
$ws.Range("G40").Value = "'5"
$ws.Range("D41").Value = "'0.007391"
$ws.Range("E41").Value = "'-1.02%"
$ws.Range("G41").Value = "'5"
$ws.Range("D42").Value = "'0.009789"
$ws.Range("E42").Value = "'-3.33%"
$ws.Range("G42").Value = "'5"
$ws.Range("D43").Value = "'0.1368"
$ws.Range("E43").Value = "'-0.20%"
$ws.Range("G43").Value = "'5"
$ws.Range("E44").Value = "'-0.08%"
$ws.Range("G44").Value = "'5"
$ws.Range("D45").Value = "'0.009377"
$ws.Range("E45").Value = "'-3.16%"
$ws.Range("G45").Value = "'5"
$ws.Range("D46").Value = "'0.00006353"
$ws.Range("E46").Value = "'0.62%"
$ws.Range("G46").Value = "'5"
$ws.Range("E47").Value = "'-0.12%"
$ws.Range("G47").Value = "'5"
$ws.Range("D48").Value = "'0.002728"
$ws.Range("G48").Value = "'5"
$ws.Range("E49").Value = "'-18.80%"
$ws.Range("G49").Value = "'5"
$ws.Range("E50").Value = "'-0.12%"
$ws.Range("G50").Value = "'5"
$ws.Range("E51").Value = "'-0.12%"
$ws.Range("G51").Value = "'5"
